$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text: MODEL_CONDITION -> MODELCONDITION
for ($i = 1; $i -le 6; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    if ($cell.Value() -eq "MODEL_CONDITION") {
        $cell.Value = "MODELCONDITION"
    }
}

# Delete the whole first column (the row-index column), shifting everything left
$ws.Columns.Item(1).Delete()
